# Insert a duplicate of row 4 as a new row 5 on Sheet1, shifting rows 5-9 down to 6-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(5).Insert()

for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(5, $col).Value2 = $ws.Cells.Item(4, $col).Value2
}
